$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed values.
# NumberFormat is set to Text ("@") before assigning so Excel keeps the values as
# literal strings (e.g. "-1.16%" and "0.06333") instead of auto-converting them to
# numbers / percentages, matching the inline-string cells used in the source sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "274.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.16%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.35%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.51%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06333"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.26%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.868"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.77%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.321"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.81%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.267"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "34.66%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8695"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.28%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "17.66%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05034"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.72%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07378"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.57%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02964"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.94%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09030"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.32%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001571"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.96%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006334"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.33%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005808"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.41%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.452"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.10%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.295"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.39%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.83%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.918"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.62%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04353"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.01%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001183"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.46%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.84%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001202"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.09%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001688"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.14%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04095"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.31%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006766"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.12%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.09%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002163"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.28%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01079"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-16.41%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005306"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.99%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.02102"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-29.63%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.490"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-37.35%"
